$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status text for the 5d5d7fd1 row changed from "Ready for handoff" to
# "Handback transform failed" everywhere it appears (Overview + both
# language sheets all shared the same string).
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# New Error Detail messages for the 5d5d7fd1 row on each language sheet.
$zhcn.Range("P3").Value = "Handback file name: avxjfxhc.4jo is different with handoff file name: 5d5d7fd1-167d-4298-aa47-c48309a84bdb.ad5e364d816fce4de7ad2473730c4f4903c0a200.zh-cn."
$dede.Range("P3").Value = "Handback file name: avxjfxhc.4jo is different with handoff file name: 5d5d7fd1-167d-4298-aa47-c48309a84bdb.ad5e364d816fce4de7ad2473730c4f4903c0a200.de-de."

# Widen the Error Detail column on both sheets so the message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
